# daily auto push: insert one new data row (2026/01/29) at row 717,
# shifting the existing rows 717-758 down to 718-759.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 717 (everything below shifts down by one).
$ws.Rows(717).Insert()

# The new row's date text ("2026/01/29") looks like a date, so Excel would
# normally auto-convert it to a date serial. Force the cell to Text format
# before assigning, then drop back to the workbook's default "Normal" style
# so no stray formatting is left behind on the cell.
$ws.Cells.Item(717, 1).NumberFormat = "@"
$ws.Cells.Item(717, 1).Value = "2026/01/29"
$ws.Cells.Item(717, 1).Style = "Normal"

$ws.Cells.Item(717, 2).Value = "木"
$ws.Cells.Item(717, 3).Value = 12
$ws.Cells.Item(717, 4).Value = 20
